# Update the LR-pairs table with new TPM-derived values (F2-Thbd)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> ECs)
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.3845463333333334
$ws.Range("H2").Value = 1.153639
$ws.Range("I2").Value = 0.1984850200147207
$ws.Range("J2").Value = 0.1984850200147207
$ws.Range("M2").Value = 29.17403400000001
$ws.Range("N2").Value = 87.52210200000002
$ws.Range("O2").Value = 0.3835306213274714
$ws.Range("P2").Value = 0.3835306213274714
$ws.Range("Q2").Value = 11.218767803242
$ws.Range("R2").Value = 100.968910229178
$ws.Range("S2").Value = 0.07612508305044144
$ws.Range("T2").Value = 0.07612508305044142

# Row 3 (ECs -> FAPs)
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.3845463333333334
$ws.Range("H3").Value = 1.153639
$ws.Range("I3").Value = 0.1984850200147207
$ws.Range("J3").Value = 0.1984850200147207
$ws.Range("O3").Value = 0.5274816184042599
$ws.Range("P3").Value = 0.5274816184042599
$ws.Range("Q3").Value = 15.42952105590278
$ws.Range("R3").Value = 138.865689503125
$ws.Range("S3").Value = 0.1046971995863668
$ws.Range("T3").Value = 0.1046971995863668

# Row 4 (ECs -> MuSCs)
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.3845463333333334
$ws.Range("H4").Value = 1.153639
$ws.Range("I4").Value = 0.1984850200147207
$ws.Range("J4").Value = 0.1984850200147207
$ws.Range("M4").Value = 6.769034333333334
$ws.Range("N4").Value = 20.307103
$ws.Range("O4").Value = 0.08898776026826866
$ws.Range("P4").Value = 0.08898776026826867
$ws.Range("Q4").Value = 2.603007333090778
$ws.Range("R4").Value = 23.427065997817
$ws.Range("S4").Value = 0.01766273737791247
$ws.Range("T4").Value = 0.01766273737791248

# Row 5 (FAPs -> ECs)
$ws.Range("I5").Value = 0.5733580031870772
$ws.Range("J5").Value = 0.5733580031870772
$ws.Range("M5").Value = 29.17403400000001
$ws.Range("N5").Value = 87.52210200000002
$ws.Range("O5").Value = 0.3835306213274714
$ws.Range("P5").Value = 0.3835306213274714
$ws.Range("Q5").Value = 32.40733384015201
$ws.Range("R5").Value = 291.6660045613681
$ws.Range("S5").Value = 0.219900351205418
$ws.Range("T5").Value = 0.219900351205418

# Row 6 (FAPs -> FAPs)
$ws.Range("I6").Value = 0.5733580031870772
$ws.Range("J6").Value = 0.5733580031870772
$ws.Range("O6").Value = 0.5274816184042599
$ws.Range("P6").Value = 0.5274816184042599
$ws.Range("S6").Value = 0.3024358074461543
$ws.Range("T6").Value = 0.3024358074461543

# Row 7 (FAPs -> MuSCs)
$ws.Range("I7").Value = 0.5733580031870772
$ws.Range("J7").Value = 0.5733580031870772
$ws.Range("M7").Value = 6.769034333333334
$ws.Range("N7").Value = 20.307103
$ws.Range("O7").Value = 0.08898776026826866
$ws.Range("P7").Value = 0.08898776026826867
$ws.Range("Q7").Value = 7.519232870428
$ws.Range("R7").Value = 67.673095833852
$ws.Range("S7").Value = 0.05102184453550484
$ws.Range("T7").Value = 0.05102184453550485

# Row 8 (MuSCs -> ECs)
$ws.Range("G8").Value = 0.442033
$ws.Range("H8").Value = 1.326099
$ws.Range("I8").Value = 0.2281569767982021
$ws.Range("J8").Value = 0.2281569767982021
$ws.Range("M8").Value = 29.17403400000001
$ws.Range("N8").Value = 87.52210200000002
$ws.Range("O8").Value = 0.3835306213274714
$ws.Range("P8").Value = 0.3835306213274714
$ws.Range("Q8").Value = 12.895885771122
$ws.Range("R8").Value = 116.062971940098
$ws.Range("S8").Value = 0.08750518707161194
$ws.Range("T8").Value = 0.08750518707161192

# Row 9 (MuSCs -> FAPs)
$ws.Range("G9").Value = 0.442033
$ws.Range("H9").Value = 1.326099
$ws.Range("I9").Value = 0.2281569767982021
$ws.Range("J9").Value = 0.2281569767982021
$ws.Range("O9").Value = 0.5274816184042599
$ws.Range("P9").Value = 0.5274816184042599
$ws.Range("Q9").Value = 17.73611367395833
$ws.Range("R9").Value = 159.625023065625
$ws.Range("S9").Value = 0.1203486113717388
$ws.Range("T9").Value = 0.1203486113717388

# Row 10 (MuSCs -> MuSCs)
$ws.Range("G10").Value = 0.442033
$ws.Range("H10").Value = 1.326099
$ws.Range("I10").Value = 0.2281569767982021
$ws.Range("J10").Value = 0.2281569767982021
$ws.Range("M10").Value = 6.769034333333334
$ws.Range("N10").Value = 20.307103
$ws.Range("O10").Value = 0.08898776026826866
$ws.Range("P10").Value = 0.08898776026826867
$ws.Range("Q10").Value = 2.992136553466333
$ws.Range("R10").Value = 26.929228981197
$ws.Range("S10").Value = 0.02030317835485134
$ws.Range("T10").Value = 0.02030317835485134
